$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$arrB = New-Object 'object[,]' 24,1
$arrB[0,0] = 12.0211950682306
$arrB[1,0] = 11.41421265077945
$arrB[2,0] = 11.02482872416699
$arrB[3,0] = 10.86212375657088
$arrB[4,0] = 10.83486900656227
$arrB[5,0] = 11.02265049019133
$arrB[6,0] = 11.81546556774948
$arrB[7,0] = 13.23162323562787
$arrB[8,0] = 14.18108771858903
$arrB[9,0] = 14.59227502865557
$arrB[10,0] = 14.74493738077927
$arrB[11,0] = 14.71219515105894
$arrB[12,0] = 14.60489601126721
$arrB[13,0] = 14.53877389232009
$arrB[14,0] = 14.15379273074599
$arrB[15,0] = 13.91226001842344
$arrB[16,0] = 13.77138866079081
$arrB[17,0] = 13.72335962574839
$arrB[18,0] = 13.93817372007882
$arrB[19,0] = 14.63649547246876
$arrB[20,0] = 15.07511688607419
$arrB[21,0] = 14.84266073509987
$arrB[22,0] = 13.92646439428854
$arrB[23,0] = 12.86410103496968
$ws.Range("B2:B25").Value = $arrB

$arrC = New-Object 'object[,]' 24,1
$arrC[0,0] = 8.742730876259355
$arrC[1,0] = 8.373186825234685
$arrC[2,0] = 8.136623512939469
$arrC[3,0] = 8.037889095241226
$arrC[4,0] = 8.021356510661741
$arrC[5,0] = 8.135301254324657
$arrC[6,0] = 8.617371168452594
$arrC[7,0] = 9.482462999675791
$arrC[8,0] = 10.06514070813247
$arrC[9,0] = 10.31808671692883
$arrC[10,0] = 10.41208581666839
$arrC[11,0] = 10.39192148404883
$arrC[12,0] = 10.32585609955413
$arrC[13,0] = 10.28515537226978
$arrC[14,0] = 10.048362368739
$arrC[15,0] = 9.899960674053238
$arrC[16,0] = 9.813465546307592
$arrC[17,0] = 9.783985711399426
$arrC[18,0] = 9.915876446830826
$arrC[19,0] = 10.34530989028385
$arrC[20,0] = 10.61554351141632
$arrC[21,0] = 10.47228133973293
$arrC[22,0] = 9.908684588762904
$arrC[23,0] = 9.257475907768933
$ws.Range("C2:C25").Value = $arrC

$arrE = New-Object 'object[,]' 24,1
$arrE[0,0] = 25.46734856096818
$arrE[1,0] = 25.206462969729
$arrE[2,0] = 25.05083115758928
$arrE[3,0] = 24.98861913245534
$arrE[4,0] = 24.97836373552698
$arrE[5,0] = 25.04998716540023
$arrE[6,0] = 25.37648956037273
$arrE[7,0] = 26.04995072090492
$arrE[8,0] = 26.56092595244839
$arrE[9,0] = 26.79601278632773
$arrE[10,0] = 26.88533938543111
$arrE[11,0] = 26.86608896836852
$arrE[12,0] = 26.80335607187305
$arrE[13,0] = 26.76496775614908
$arrE[14,0] = 26.54560916718747
$arrE[15,0] = 26.41166254200181
$arrE[16,0] = 26.33487314805448
$arrE[17,0] = 26.30891936553661
$arrE[18,0] = 26.42589573574595
$arrE[19,0] = 26.82177458544584
$arrE[20,0] = 27.08224339587867
$arrE[21,0] = 26.94309200682242
$arrE[22,0] = 26.41946022304846
$arrE[23,0] = 25.86462150431095
$ws.Range("E2:E25").Value = $arrE

$arrF = New-Object 'object[,]' 24,1
$arrF[0,0] = 38.11783609816463
$arrF[1,0] = 37.87459206406198
$arrF[2,0] = 37.73586233072821
$arrF[3,0] = 37.68204921883589
$arrF[4,0] = 37.67327923316896
$arrF[5,0] = 37.73512551250455
$arrF[6,0] = 38.03179172099779
$arrF[7,0] = 38.69544994742939
$arrF[8,0] = 39.22945101021862
$arrF[9,0] = 39.48162926530315
$arrF[10,0] = 39.57837935404896
$arrF[11,0] = 39.55748787814323
$arrF[12,0] = 39.48956423084078
$arrF[13,0] = 39.44812026823086
$arrF[14,0] = 39.21315046212894
$arrF[15,0] = 39.07132241996121
$arrF[16,0] = 38.99062285173898
$arrF[17,0] = 38.96345211028272
$arrF[18,0] = 39.08633011674699
$arrF[19,0] = 39.5094815904118
$arrF[20,0] = 39.79331765297349
$arrF[21,0] = 39.64118833834262
$arrF[22,0] = 39.07954251696066
$arrF[23,0] = 38.5074917367992
$ws.Range("F2:F25").Value = $arrF

$arrG = New-Object 'object[,]' 24,1
$arrG[0,0] = 19.27444935623251
$arrG[1,0] = 19.2447527837198
$arrG[2,0] = 19.23634801628086
$arrG[3,0] = 19.23538984315835
$arrG[4,0] = 19.23537951460641
$arrG[5,0] = 19.23632511511948
$arrG[6,0] = 19.26216830001613
$arrG[7,0] = 19.39085554998015
$arrG[8,0] = 19.53270773539456
$arrG[9,0] = 19.60737282778421
$arrG[10,0] = 19.63708630398533
$arrG[11,0] = 19.63062328523638
$arrG[12,0] = 19.60978864239323
$arrG[13,0] = 19.5972136751978
$arrG[14,0] = 19.52803069645059
$arrG[15,0] = 19.48817356922143
$arrG[16,0] = 19.46620416989439
$arrG[17,0] = 19.45893028152374
$arrG[18,0] = 19.49231767040139
$arrG[19,0] = 19.6158693796475
$arrG[20,0] = 19.70499657657295
$arrG[21,0] = 19.65666799473271
$arrG[22,0] = 19.49044117731109
$arrG[23,0] = 19.34770062484287
$ws.Range("G2:G25").Value = $arrG

$arrH = New-Object 'object[,]' 24,1
$arrH[0,0] = 11.79075454630943
$arrH[1,0] = 11.83646508628791
$arrH[2,0] = 11.86689144278723
$arrH[3,0] = 11.87988293370328
$arrH[4,0] = 11.8820759199695
$arrH[5,0] = 11.86706425277203
$arrH[6,0] = 11.80602519616887
$arrH[7,0] = 11.70509324294769
$arrH[8,0] = 11.64243401179867
$arrH[9,0] = 11.61643678749895
$arrH[10,0] = 11.60695380469098
$arrH[11,0] = 11.6089800344227
$arrH[12,0] = 11.61564936329994
$arrH[13,0] = 11.61978164450803
$arrH[14,0] = 11.64418353102597
$arrH[15,0] = 11.65979614914463
$arrH[16,0] = 11.66901199620641
$arrH[17,0] = 11.67217280212761
$arrH[18,0] = 11.65810973816327
$arrH[19,0] = 11.61368059698651
$arrH[20,0] = 11.58675179757873
$arrH[21,0] = 11.60093093090456
$arrH[22,0] = 11.65887141721755
$arrH[23,0] = 11.73038330051161
$ws.Range("H2:H25").Value = $arrH

$arrJ = New-Object 'object[,]' 24,1
$arrJ[0,0] = 7.262919677881048
$arrJ[1,0] = 7.28470662220512
$arrJ[2,0] = 7.298927512804203
$arrJ[3,0] = 7.304935119313762
$arrJ[4,0] = 7.305945521129178
$arrJ[5,0] = 7.299007672710832
$arrJ[6,0] = 7.270256929525271
$arrJ[7,0] = 7.220555024360956
$arrJ[8,0] = 7.188088746501712
$arrJ[9,0] = 7.174193877057254
$arrJ[10,0] = 7.169057638048767
$arrJ[11,0] = 7.17015824386609
$arrJ[12,0] = 7.173768802929817
$arrJ[13,0] = 7.175996703759706
$arrJ[14,0] = 7.189014377475765
$arrJ[15,0] = 7.197224013884101
$arrJ[16,0] = 7.20202827903908
$arrJ[17,0] = 7.203669063800181
$arrJ[18,0] = 7.196341568316862
$arrJ[19,0] = 7.172704891690975
$arrJ[20,0] = 7.157988038128115
$arrJ[21,0] = 7.165775890973163
$arrJ[22,0] = 7.196740258927099
$arrJ[23,0] = 7.23328800592498
$ws.Range("J2:J25").Value = $arrJ

$arrN = New-Object 'object[,]' 24,1
$arrN[0,0] = 15.74979979977724
$arrN[1,0] = 15.76688132243443
$arrN[2,0] = 15.7791339398228
$arrN[3,0] = 15.78457159872191
$arrN[4,0] = 15.78550139801746
$arrN[5,0] = 15.77920547256249
$arrN[6,0] = 15.75532384848838
$arrN[7,0] = 15.72245191608583
$arrN[8,0] = 15.70675641161005
$arrN[9,0] = 15.70143986090769
$arrN[10,0] = 15.69968771098363
$arrN[11,0] = 15.70005347260853
$arrN[12,0] = 15.70129048380144
$arrN[13,0] = 15.70208216029551
$arrN[14,0] = 15.70714046002784
$arrN[15,0] = 15.71070981906731
$arrN[16,0] = 15.71293451615005
$arrN[17,0] = 15.71371728194332
$arrN[18,0] = 15.7103120925667
$arrN[19,0] = 15.70092006602255
$arrN[20,0] = 15.69630326049406
$arrN[21,0] = 15.69862849751892
$arrN[22,0] = 15.71049136691102
$arrN[23,0] = 15.72985565679904
$ws.Range("N2:N25").Value = $arrN

$arrO = New-Object 'object[,]' 24,1
$arrO[0,0] = 16.68235786430008
$arrO[1,0] = 16.74299550131385
$arrO[2,0] = 16.78515260777696
$arrO[3,0] = 16.80356481350271
$arrO[4,0] = 16.8066964367401
$arrO[5,0] = 16.78539593717825
$arrO[6,0] = 16.70224011571007
$arrO[7,0] = 16.5784966185082
$arrO[8,0] = 16.51187517261591
$arrO[9,0] = 16.48690521691458
$arrO[10,0] = 16.4782218125834
$arrO[11,0] = 16.48005752348903
$arrO[12,0] = 16.48617532616156
$arrO[13,0] = 16.49002335376901
$arrO[14,0] = 16.51361481906547
$arrO[15,0] = 16.52945773524208
$arrO[16,0] = 16.53907218904782
$arrO[17,0] = 16.54241355009137
$arrO[18,0] = 16.52771923598825
$arrO[19,0] = 16.48435738599849
$arrO[20,0] = 16.46052072994604
$arrO[21,0] = 16.47282929544551
$arrO[22,0] = 16.528503635209
$arrO[23,0] = 16.60772849466944
$ws.Range("O2:O25").Value = $arrO
